$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A10").Value = 111897633
$ws.Range("B10").Value = 95674
$ws.Range("E10").Value = 222741
$ws.Range("Q10").Value = 489384.8847604021
$ws.Range("R10").Value = 7032364.981337981
$ws.Range("S10").Value = 10

# Text cells
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "LC"
$ws.Range("F10").Value = "Finbräken"
$ws.Range("G10").Value = "Cystopteris montana"
$ws.Range("H10").Value = "(Lam.) Desv."
$ws.Range("P10").Value = "Husås, Jmt"
$ws.Range("T10").Value = "Jämtland"
$ws.Range("U10").Value = "Östersund"
$ws.Range("V10").Value = "Jämtland"
$ws.Range("W10").Value = "Lit"
$ws.Range("AW10").Value = "Christer Pålsson"
$ws.Range("AX10").Value = "Christer Pålsson"

# Date/time-like text cells -- must stay literal text, not auto-converted to date serials.
# Pre-format as Text so Excel keeps them literal, then restore the Normal style so no
# explicit cell style (numFmt) lingers in the saved file.
$ws.Range("Y10:AB10").NumberFormat = "@"
$ws.Range("Y10").Value = "2023-06-22"
$ws.Range("Z10").Value = "00:00"
$ws.Range("AA10").Value = "2023-06-22"
$ws.Range("AB10").Value = "00:00"
$ws.Range("Y10:AB10").Style = "Normal"

# Boolean cells
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false

# Blank-but-present cells (mirrors the row above, which has the same empty columns)
$ws.Range("I9").Copy($ws.Range("I10"))
$ws.Range("J9").Copy($ws.Range("J10"))
$ws.Range("K9").Copy($ws.Range("K10"))
$ws.Range("L9").Copy($ws.Range("L10"))
$ws.Range("N9").Copy($ws.Range("N10"))
$ws.Range("AF9").Copy($ws.Range("AF10"))
$ws.Range("AT9").Copy($ws.Range("AT10"))
$ws.Range("AY9").Copy($ws.Range("AY10"))

$excel.CutCopyMode = $false
